# Weekly update: insert one new "Perejil" price record as a new row 661
# in the consolidated sheet, shifting the existing rows 661..710 down to
# 662..711 (dimension grows from A1:R710 to A1:R711).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 661 (pushes 661-710 -> 662-711).
$ws.Rows.Item(661).Insert()

# Populate the newly inserted row 661 with the new weekly record.
$ws.Range("A661").Value = 6
$ws.Range("B661").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C661").Value = "Metropolitana"
$ws.Range("D661").Value = 45013
$ws.Range("E661").Value = 13
$ws.Range("F661").Value = 100112044
$ws.Range("G661").Value = "Perejil"
$ws.Range("H661").Value = "Sin especificar"
$ws.Range("I661").Value = "Primera"
$ws.Range("J661").Value = 280
$ws.Range("K661").Value = 12000
$ws.Range("L661").Value = 13000
$ws.Range("M661").Value = 12536
$ws.Range("N661").Value = "$/docena de atados"
$ws.Range("O661").Value = "Región Metropolitana"
$ws.Range("P661").Value = 4179
$ws.Range("Q661").Value = 3
$ws.Range("R661").Value = "Hortaliza"
